$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to text format so numeric-looking values
# (e.g. "299.74") are preserved as strings rather than being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '42.204.43'
$ws.Range('E2').Value = '  -1.24%  '
$ws.Range('D3').Value = '2.268.50'
$ws.Range('E3').Value = '  -1.57%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '299.74'
$ws.Range('E5').Value = '  -1.40%  '
$ws.Range('D6').Value = '96.22'
$ws.Range('E6').Value = '  -2.92%  '
$ws.Range('D7').Value = '0.495'
$ws.Range('E7').Value = '  -2.18%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '0.492'
$ws.Range('E9').Value = '  -1.93%  '
$ws.Range('D10').Value = '33.24'
$ws.Range('E10').Value = '  -3.31%  '
$ws.Range('D11').Value = '0.0787'
$ws.Range('E11').Value = '  -0.23%  '
$ws.Range('D12').Value = '48.20'
$ws.Range('E12').Value = '  -6.88%  '
$ws.Range('E13').Value = '  +0.33%  '
$ws.Range('D14').Value = '6.66'
$ws.Range('E14').Value = '  -1.18%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.626.81'
$ws.Range('E15').Value = '  -1.43%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').Value = '15.54'
$ws.Range('E16').Value = '  -0.76%  '
$ws.Range('D17').Value = '2.271.29'
$ws.Range('E17').Value = '  -1.06%  '
$ws.Range('D18').Value = '0.784'
$ws.Range('E18').Value = '  -3.06%  '
$ws.Range('D19').Value = '42.148.54'
$ws.Range('E19').Value = '  -1.22%  '
$ws.Range('D20').Value = '11.68'
$ws.Range('E20').Value = '  +1.18%  '
$ws.Range('D21').Value = '0.0₃0888'
$ws.Range('E21').Value = '  -1.37%  '
$ws.Range('D22').Value = '5.99'
$ws.Range('E22').Value = '  -1.51%  '
$ws.Range('D23').Value = '66.41'
$ws.Range('E23').Value = '  -4.15%  '
$ws.Range('D24').Value = '234.80'
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').Value = '1.97'
$ws.Range('E25').Value = '  -0.59%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').Value = '2.46'
$ws.Range('E26').Value = '  -2.45%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('D28').Value = '23.96'
$ws.Range('E28').Value = '  -4.89%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '2.17'
$ws.Range('E29').Value = '  -5.75%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').Value = '168.45'
$ws.Range('E30').Value = '  +3.50%  '
$ws.Range('D31').Value = '9.19'
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('D32').Value = '33.70'
$ws.Range('E32').Value = '  -2.62%  '
$ws.Range('E33').Value = '  -0.12%  '
$ws.Range('D34').Value = '4.90'
$ws.Range('E34').Value = '  -2.66%  '
$ws.Range('D35').Value = '4.57'
$ws.Range('E35').Value = '  -0.95%  '
$ws.Range('D36').Value = '16.54'
$ws.Range('E36').Value = '  -2.42%  '
$ws.Range('E37').Value = '  -4.86%  '
$ws.Range('D38').Value = '0.0686'
$ws.Range('E38').Value = '  -3.85%  '
$ws.Range('D39').Value = '2.78'
$ws.Range('E39').Value = '  -3.66%  '
$ws.Range('D40').Value = '0.0988'
$ws.Range('E40').Value = '  -1.73%  '
$ws.Range('E41').Value = '  -2.54%  '
$ws.Range('E42').Value = '  -4.49%  '
$ws.Range('D43').Value = '2.43'
$ws.Range('E43').Value = '  -1.43%  '
$ws.Range('D44').Value = '1.972.54'
$ws.Range('E44').Value = '  -0.85%  '
$ws.Range('D45').Value = '0.0277'
$ws.Range('E45').Value = '  -1.02%  '
$ws.Range('D46').Value = '17.40'
$ws.Range('E46').Value = '  -7.34%  '
$ws.Range('D47').Value = '9.56'
$ws.Range('E47').Value = '  -6.70%  '
$ws.Range('D48').Value = '2.78'
$ws.Range('E48').Value = '  -3.12%  '
$ws.Range('D49').Value = '2.496.81'
$ws.Range('E49').Value = '  -1.39%  '
$ws.Range('D50').Value = '52.30'
$ws.Range('E50').Value = '  -5.52%  '
$ws.Range('D51').Value = '1.47'
$ws.Range('E51').Value = '  -1.23%  '

# Restore column D to its original (default/Normal) style now that values are set.
$ws.Range("D2:D51").Style = "Normal"
